# Update Release-Notes.xlsx - Folder inventory updated on Mon Jun 16 14:10:51 UTC 2025
$wb = $excel.ActiveWorkbook

# --- Folder Inventory sheet: new row inserted at top, list re-sorted/shifted ---
$ws = $wb.Worksheets.Item("Folder Inventory")

$ws.Cells.Item(2, 1).Value = 'GitHub Copilot Innovation Workshop'
$ws.Cells.Item(2, 2).Value = 'GitHub Copilot Innovation Workshop'
$ws.Cells.Item(2, 3).Value = '2025-06-16 19:40:23 +0530'
$ws.Cells.Item(2, 4).Value = 1
$ws.Cells.Item(2, 5).Value = 'Root'

$ws.Cells.Item(3, 1).Value = 'Microsoft Defender for Cloud - v1'
$ws.Cells.Item(3, 2).Value = 'Microsoft Defender for Cloud - v1'
$ws.Cells.Item(3, 3).Value = '2025-06-16 17:12:47 +0530'
$ws.Cells.Item(3, 4).Value = 1
$ws.Cells.Item(3, 5).Value = 'Root'

$ws.Cells.Item(4, 1).Value = 'Build Intelligent Apps with Microsoft''s Copilot stack & Azure OpenAI  '
$ws.Cells.Item(4, 2).Value = 'Build Intelligent Apps with Microsoft''s Copilot stack & Azure OpenAI  '
$ws.Cells.Item(4, 3).Value = '2025-06-16 16:14:06 +0530'
$ws.Cells.Item(4, 4).Value = 1
$ws.Cells.Item(4, 5).Value = 'Root'

$ws.Cells.Item(5, 1).Value = 'Hackathon - Intelligent App Development with Microsoft CoPilot Stack'
$ws.Cells.Item(5, 2).Value = 'Hackathon - Intelligent App Development with Microsoft CoPilot Stack'
$ws.Cells.Item(5, 3).Value = '2025-06-16 14:04:22 +0530'
$ws.Cells.Item(5, 4).Value = 1
$ws.Cells.Item(5, 5).Value = 'Root'

$ws.Cells.Item(6, 1).Value = 'Azure_Well-Architected_Resiliency_Gaps_Remediation'
$ws.Cells.Item(6, 2).Value = 'Azure_Well-Architected_Resiliency_Gaps_Remediation'
$ws.Cells.Item(6, 3).Value = '2025-06-13 17:35:45 +0530'
$ws.Cells.Item(6, 4).Value = 1
$ws.Cells.Item(6, 5).Value = 'Root'

$ws.Cells.Item(7, 1).Value = 'Implement CI-CD with GitHub Actions'
$ws.Cells.Item(7, 2).Value = 'Implement CI-CD with GitHub Actions'
$ws.Cells.Item(7, 3).Value = '2025-06-13 15:19:07 +0000'
$ws.Cells.Item(7, 4).Value = 1
$ws.Cells.Item(7, 5).Value = 'Root'

$ws.Cells.Item(8, 1).Value = 'Developing_a_Custom_RAG_App_Using_Azure_AI_Foundry'
$ws.Cells.Item(8, 2).Value = 'Developing_a_Custom_RAG_App_Using_Azure_AI_Foundry'
$ws.Cells.Item(8, 3).Value = '2025-06-13 13:09:40 +0530'
$ws.Cells.Item(8, 4).Value = 1
$ws.Cells.Item(8, 5).Value = 'Root'

$ws.Cells.Item(9, 1).Value = 'Automated Machine Learning Using AML'
$ws.Cells.Item(9, 2).Value = 'Automated Machine Learning Using AML'
$ws.Cells.Item(9, 3).Value = '2025-06-12 21:50:14 +0530'
$ws.Cells.Item(9, 4).Value = 1
$ws.Cells.Item(9, 5).Value = 'Root'

$ws.Cells.Item(10, 1).Value = 'Create and Publish PowerBI Dashboards & Reports'
$ws.Cells.Item(10, 2).Value = 'Create and Publish PowerBI Dashboards & Reports'
$ws.Cells.Item(10, 3).Value = '2025-06-12 20:05:46 +0530'
$ws.Cells.Item(10, 4).Value = 1
$ws.Cells.Item(10, 5).Value = 'Root'

$ws.Cells.Item(11, 1).Value = 'Azure Virtual Machine And Compute'
$ws.Cells.Item(11, 2).Value = 'Azure Virtual Machine And Compute'
$ws.Cells.Item(11, 3).Value = '2025-06-12 17:37:08 +0530'
$ws.Cells.Item(11, 4).Value = 1
$ws.Cells.Item(11, 5).Value = 'Root'

$ws.Cells.Item(12, 1).Value = 'Work with Data Lake and Data Factory Pipelines in Microsoft Fabric​'
$ws.Cells.Item(12, 2).Value = 'Work with Data Lake and Data Factory Pipelines in Microsoft Fabric​'
$ws.Cells.Item(12, 3).Value = '2025-06-12 17:26:19 +0530'
$ws.Cells.Item(12, 4).Value = 1
$ws.Cells.Item(12, 5).Value = 'Root'

$ws.Cells.Item(13, 1).Value = 'Get Started with Microsoft Fabric with Its Lakehouses'
$ws.Cells.Item(13, 2).Value = 'Get Started with Microsoft Fabric with Its Lakehouses'
$ws.Cells.Item(13, 3).Value = '2025-06-12 16:16:30 +0530'
$ws.Cells.Item(13, 4).Value = 1
$ws.Cells.Item(13, 5).Value = 'Root'

$ws.Cells.Item(14, 1).Value = 'Build A Fabric Real-Time Intelligence Solution in a Day'
$ws.Cells.Item(14, 2).Value = 'Build A Fabric Real-Time Intelligence Solution in a Day'
$ws.Cells.Item(14, 3).Value = '2025-06-12 15:59:35 +0530'
$ws.Cells.Item(14, 4).Value = 1
$ws.Cells.Item(14, 5).Value = 'Root'

$ws.Cells.Item(15, 1).Value = 'Azure_AI_Foundry_and_Semantic_Kernel_Fundamentals'
$ws.Cells.Item(15, 2).Value = 'Azure_AI_Foundry_and_Semantic_Kernel_Fundamentals'
$ws.Cells.Item(15, 3).Value = '2025-06-12 15:19:27 +0530'
$ws.Cells.Item(15, 4).Value = 1
$ws.Cells.Item(15, 5).Value = 'Root'

$ws.Cells.Item(16, 1).Value = 'Enterprise-Class Networking in Azure'
$ws.Cells.Item(16, 2).Value = 'Enterprise-Class Networking in Azure'
$ws.Cells.Item(16, 3).Value = '2025-06-12 12:35:48 +0530'
$ws.Cells.Item(16, 4).Value = 1
$ws.Cells.Item(16, 5).Value = 'Root'

$ws.Cells.Item(17, 1).Value = 'Cloud-Native Applications'
$ws.Cells.Item(17, 2).Value = 'Cloud-Native Applications'
$ws.Cells.Item(17, 3).Value = '2025-06-12 12:18:28 +0530'
$ws.Cells.Item(17, 4).Value = 1
$ws.Cells.Item(17, 5).Value = 'Root'

$ws.Cells.Item(18, 1).Value = 'Hackathon - Activate GenAI with Azure'
$ws.Cells.Item(18, 2).Value = 'Hackathon - Activate GenAI with Azure'
$ws.Cells.Item(18, 3).Value = '2025-06-11 22:47:04 +0530'
$ws.Cells.Item(18, 4).Value = 2
$ws.Cells.Item(18, 5).Value = 'Root'

$ws.Cells.Item(19, 1).Value = 'Azure Landing Zone'
$ws.Cells.Item(19, 2).Value = 'Azure Landing Zone'
$ws.Cells.Item(19, 3).Value = '2025-06-11 20:16:49 +0530'
$ws.Cells.Item(19, 4).Value = 1
$ws.Cells.Item(19, 5).Value = 'Root'

$ws.Cells.Item(20, 1).Value = 'Microsoft Azure AI Agents'
$ws.Cells.Item(20, 2).Value = 'Microsoft Azure AI Agents'
$ws.Cells.Item(20, 3).Value = '2025-06-11 20:13:48 +0530'
$ws.Cells.Item(20, 4).Value = 1
$ws.Cells.Item(20, 5).Value = 'Root'

$ws.Cells.Item(21, 1).Value = 'Azure Local Hands-on Lab'
$ws.Cells.Item(21, 2).Value = 'Azure Local Hands-on Lab'
$ws.Cells.Item(21, 3).Value = '2025-06-11 19:56:28 +0530'
$ws.Cells.Item(21, 4).Value = 1
$ws.Cells.Item(21, 5).Value = 'Root'

$ws.Cells.Item(22, 1).Value = 'Get data into Fabric Lakehouse'
$ws.Cells.Item(22, 2).Value = 'Get data into Fabric Lakehouse'
$ws.Cells.Item(22, 3).Value = '2025-06-11 15:00:50 +0000'
$ws.Cells.Item(22, 4).Value = 1
$ws.Cells.Item(22, 5).Value = 'Root'

$ws.Cells.Item(23, 1).Value = 'Low Code for Pro-Dev in a Day'
$ws.Cells.Item(23, 2).Value = 'Low Code for Pro-Dev in a Day'
$ws.Cells.Item(23, 3).Value = '2025-06-11 00:35:20 +0530'
$ws.Cells.Item(23, 4).Value = 1
$ws.Cells.Item(23, 5).Value = 'Root'

$ws.Cells.Item(24, 1).Value = 'Developing AI Applications with Azure AI Foundry'
$ws.Cells.Item(24, 2).Value = 'Developing AI Applications with Azure AI Foundry'
$ws.Cells.Item(24, 3).Value = '2025-06-11 00:33:06 +0530'
$ws.Cells.Item(24, 4).Value = 1
$ws.Cells.Item(24, 5).Value = 'Root'

$ws.Cells.Item(25, 1).Value = 'Develop Generative AI solutions with Azure OpenAI Service'
$ws.Cells.Item(25, 2).Value = 'Develop Generative AI solutions with Azure OpenAI Service'
$ws.Cells.Item(25, 3).Value = '2025-06-10 23:22:30 +0530'
$ws.Cells.Item(25, 4).Value = 1
$ws.Cells.Item(25, 5).Value = 'Root'

$ws.Cells.Item(26, 1).Value = 'Advanced Workflow Automation with GitHub Actions '
$ws.Cells.Item(26, 2).Value = 'Advanced Workflow Automation with GitHub Actions '
$ws.Cells.Item(26, 3).Value = '2025-06-10 23:10:36 +0530'
$ws.Cells.Item(26, 4).Value = 1
$ws.Cells.Item(26, 5).Value = 'Root'

$ws.Cells.Item(27, 1).Value = 'Get Started With OpenAI And Build Natural Language Solution'
$ws.Cells.Item(27, 2).Value = 'Get Started With OpenAI And Build Natural Language Solution'
$ws.Cells.Item(27, 3).Value = '2025-06-10 22:51:47 +0530'
$ws.Cells.Item(27, 4).Value = 1
$ws.Cells.Item(27, 5).Value = 'Root'

$ws.Cells.Item(28, 1).Value = 'Lunch and Learn: Building and Evaluating Prompt Flows with Azure AI Foundry'
$ws.Cells.Item(28, 2).Value = 'Lunch and Learn: Building and Evaluating Prompt Flows with Azure AI Foundry'
$ws.Cells.Item(28, 3).Value = '2025-06-10 22:48:16 +0530'
$ws.Cells.Item(28, 4).Value = 1
$ws.Cells.Item(28, 5).Value = 'Root'

$ws.Cells.Item(29, 1).Value = 'Automate document processing by using Azure AI & OpenAI'
$ws.Cells.Item(29, 2).Value = 'Automate document processing by using Azure AI & OpenAI'
$ws.Cells.Item(29, 3).Value = '2025-06-10 07:51:57 +0530'
$ws.Cells.Item(29, 4).Value = 1
$ws.Cells.Item(29, 5).Value = 'Root'

$ws.Cells.Item(30, 1).Value = 'Azure Virtual Desktop'
$ws.Cells.Item(30, 2).Value = 'Azure Virtual Desktop'
$ws.Cells.Item(30, 3).Value = '2025-06-09 23:17:02 +0530'
$ws.Cells.Item(30, 4).Value = 1
$ws.Cells.Item(30, 5).Value = 'Root'

$ws.Cells.Item(31, 1).Value = 'Get Started with Real-Time Analytics and Data Science with Microsoft Fabric'
$ws.Cells.Item(31, 2).Value = 'Get Started with Real-Time Analytics and Data Science with Microsoft Fabric'
$ws.Cells.Item(31, 3).Value = '2025-06-09 18:18:42 +0530'
$ws.Cells.Item(31, 4).Value = 1
$ws.Cells.Item(31, 5).Value = 'Root'

$ws.Cells.Item(32, 1).Value = 'Code Suggestions with GitHub Copilot in Codespace using VS Code'
$ws.Cells.Item(32, 2).Value = 'Code Suggestions with GitHub Copilot in Codespace using VS Code'
$ws.Cells.Item(32, 3).Value = '2025-06-09 17:15:52 +0530'
$ws.Cells.Item(32, 4).Value = 1
$ws.Cells.Item(32, 5).Value = 'Root'

$ws.Cells.Item(33, 1).Value = 'Customer Support Conversation Summarization with Azure OpenAI'
$ws.Cells.Item(33, 2).Value = 'Customer Support Conversation Summarization with Azure OpenAI'
$ws.Cells.Item(33, 3).Value = '2025-06-09 16:24:38 +0530'
$ws.Cells.Item(33, 4).Value = 1
$ws.Cells.Item(33, 5).Value = 'Root'

$ws.Cells.Item(34, 1).Value = 'Secure Windows Servers Azure Arc & Microsoft Defender'
$ws.Cells.Item(34, 2).Value = 'Secure Windows Servers Azure Arc & Microsoft Defender'
$ws.Cells.Item(34, 3).Value = '2025-06-09 14:20:04 +0000'
$ws.Cells.Item(34, 4).Value = 1
$ws.Cells.Item(34, 5).Value = 'Root'

$ws.Cells.Item(35, 1).Value = 'Build Prompt Engineering With Azure OpenAI Service'
$ws.Cells.Item(35, 2).Value = 'Build Prompt Engineering With Azure OpenAI Service'
$ws.Cells.Item(35, 3).Value = '2025-06-05 22:32:54 +0530'
$ws.Cells.Item(35, 4).Value = 1
$ws.Cells.Item(35, 5).Value = 'Root'

$ws.Cells.Item(36, 1).Value = 'Lunch and Learn: Build Custom Copilot Application using Azure AI Foundry'
$ws.Cells.Item(36, 2).Value = 'Lunch and Learn: Build Custom Copilot Application using Azure AI Foundry'
$ws.Cells.Item(36, 3).Value = '2025-06-05 21:53:26 +0530'
$ws.Cells.Item(36, 4).Value = 1
$ws.Cells.Item(36, 5).Value = 'Root'

$ws.Cells.Item(37, 1).Value = 'Microsoft Fabric with capacity-copilot-SDP'
$ws.Cells.Item(37, 2).Value = 'Microsoft Fabric with capacity-copilot-SDP'
$ws.Cells.Item(37, 3).Value = '2025-06-05 19:53:59 +0530'
$ws.Cells.Item(37, 4).Value = 1
$ws.Cells.Item(37, 5).Value = 'Root'

$ws.Cells.Item(38, 1).Value = 'Migrate-and-Modernise-SQL-Servers-to-Azure'
$ws.Cells.Item(38, 2).Value = 'Migrate-and-Modernise-SQL-Servers-to-Azure'
$ws.Cells.Item(38, 3).Value = '2025-06-05 19:52:34 +0530'
$ws.Cells.Item(38, 4).Value = 1
$ws.Cells.Item(38, 5).Value = 'Root'

$ws.Cells.Item(39, 1).Value = 'Microsoft Defender for Cloud - v2'
$ws.Cells.Item(39, 2).Value = 'Microsoft Defender for Cloud - v2'
$ws.Cells.Item(39, 3).Value = '2025-06-05 19:52:20 +0530'
$ws.Cells.Item(39, 4).Value = 1
$ws.Cells.Item(39, 5).Value = 'Root'

$ws.Cells.Item(40, 1).Value = 'Github&AgenticAI'
$ws.Cells.Item(40, 2).Value = 'Github&AgenticAI'
$ws.Cells.Item(40, 3).Value = '2025-06-04 23:32:46 +0530'
$ws.Cells.Item(40, 4).Value = 1
$ws.Cells.Item(40, 5).Value = 'Root'

$ws.Cells.Item(41, 1).Value = 'Leverage-Microsoft-365-Copilot-and-Copilot-Studio-for-Marketing'
$ws.Cells.Item(41, 2).Value = 'Leverage-Microsoft-365-Copilot-and-Copilot-Studio-for-Marketing'
$ws.Cells.Item(41, 3).Value = '2025-06-04 17:43:54 +0300'
$ws.Cells.Item(41, 4).Value = 1
$ws.Cells.Item(41, 5).Value = 'Root'

$ws.Cells.Item(42, 1).Value = 'Leverage-Microsoft-365-Copilot-and-Copilot-Studio-for-Sales'
$ws.Cells.Item(42, 2).Value = 'Leverage-Microsoft-365-Copilot-and-Copilot-Studio-for-Sales'
$ws.Cells.Item(42, 3).Value = '2025-06-04 17:41:59 +0300'
$ws.Cells.Item(42, 4).Value = 1
$ws.Cells.Item(42, 5).Value = 'Root'

$ws.Cells.Item(43, 1).Value = 'Leverage-Microsoft365-Copilot-and-Copilot-Studio-for-Human-Resources'
$ws.Cells.Item(43, 2).Value = 'Leverage-Microsoft365-Copilot-and-Copilot-Studio-for-Human-Resources'
$ws.Cells.Item(43, 3).Value = '2025-06-04 17:37:32 +0300'
$ws.Cells.Item(43, 4).Value = 1
$ws.Cells.Item(43, 5).Value = 'Root'

$ws.Cells.Item(44, 1).Value = 'Use Azure Open AI Like A Pro To Build Powerful AI Applications'
$ws.Cells.Item(44, 2).Value = 'Use Azure Open AI Like A Pro To Build Powerful AI Applications'
$ws.Cells.Item(44, 3).Value = '2025-06-04 15:08:23 +0530'
$ws.Cells.Item(44, 4).Value = 1
$ws.Cells.Item(44, 5).Value = 'Root'

$ws.Cells.Item(45, 1).Value = 'Business Automation with Azure OpenAI and Document Intelligence'
$ws.Cells.Item(45, 2).Value = 'Business Automation with Azure OpenAI and Document Intelligence'
$ws.Cells.Item(45, 3).Value = '2025-06-04 15:04:23 +0530'
$ws.Cells.Item(45, 4).Value = 1
$ws.Cells.Item(45, 5).Value = 'Root'

$ws.Cells.Item(46, 1).Value = 'Fabric Copilot'
$ws.Cells.Item(46, 2).Value = 'Fabric Copilot'
$ws.Cells.Item(46, 3).Value = '2025-06-04 02:16:22 +0530'
$ws.Cells.Item(46, 4).Value = 1
$ws.Cells.Item(46, 5).Value = 'Root'

$ws.Cells.Item(47, 1).Value = 'Integrate Azure OpenAI into your app'
$ws.Cells.Item(47, 2).Value = 'Integrate Azure OpenAI into your app'
$ws.Cells.Item(47, 3).Value = '2025-06-03 13:54:36 +0530'
$ws.Cells.Item(47, 4).Value = 1
$ws.Cells.Item(47, 5).Value = 'Root'

$ws.Cells.Item(48, 1).Value = 'Hackathon - GitHub Copilot '
$ws.Cells.Item(48, 2).Value = 'Hackathon - GitHub Copilot '
$ws.Cells.Item(48, 3).Value = '2025-05-30 20:59:44 +0530'
$ws.Cells.Item(48, 4).Value = 1
$ws.Cells.Item(48, 5).Value = 'Root'

$ws.Cells.Item(49, 1).Value = 'Call Center data analysis using Azure AI services and Azure OpenAI '
$ws.Cells.Item(49, 2).Value = 'Call Center data analysis using Azure AI services and Azure OpenAI '
$ws.Cells.Item(49, 3).Value = '2025-05-29 23:23:52 +0530'
$ws.Cells.Item(49, 4).Value = 1
$ws.Cells.Item(49, 5).Value = 'Root'

# --- Metadata sheet: refresh generated timestamp and workflow run count ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B3").Value = '2025-06-16 14:10:51 UTC'
$wsMeta.Range("B5").Value = '29'

# --- Summary sheet: refresh most recent update timestamp ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B5").Value = '2025-06-16 19:40:23 +0530'
